$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.427210000000001
$ws.Range("H2").Value = 28.28163
$ws.Range("I2").Value = 0.2188083857550241
$ws.Range("J2").Value = 0.2188083857550241
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 18.08993297466333
$ws.Range("R2").Value = 162.80939677197
$ws.Range("S2").Value = 0.0014275101015834
$ws.Range("T2").Value = 0.0014275101015834
$ws.Range("G3").Value = 9.427210000000001
$ws.Range("H3").Value = 28.28163
$ws.Range("I3").Value = 0.2188083857550241
$ws.Range("J3").Value = 0.2188083857550241
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 1709.043735661514
$ws.Range("R3").Value = 15381.39362095362
$ws.Range("S3").Value = 0.1348638051960524
$ws.Range("T3").Value = 0.1348638051960524
$ws.Range("G4").Value = 9.427210000000001
$ws.Range("H4").Value = 28.28163
$ws.Range("I4").Value = 0.2188083857550241
$ws.Range("J4").Value = 0.2188083857550241
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 273.8385070975233
$ws.Range("R4").Value = 2464.54656387771
$ws.Range("S4").Value = 0.02160910356228157
$ws.Range("T4").Value = 0.02160910356228157
$ws.Range("G5").Value = 9.427210000000001
$ws.Range("H5").Value = 28.28163
$ws.Range("I5").Value = 0.2188083857550241
$ws.Range("J5").Value = 0.2188083857550241
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 771.8481554234534
$ws.Range("R5").Value = 6946.63339881108
$ws.Range("S5").Value = 0.06090796689510673
$ws.Range("T5").Value = 0.06090796689510673
$ws.Range("I6").Value = 0.3808887290954196
$ws.Range("J6").Value = 0.3808887290954196
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 31.48988808799633
$ws.Range("R6").Value = 283.408992791967
$ws.Range("S6").Value = 0.00248492536740215
$ws.Range("T6").Value = 0.00248492536740215
$ws.Range("I7").Value = 0.3808887290954196
$ws.Range("J7").Value = 0.3808887290954196
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.2347629556556754
$ws.Range("T7").Value = 0.2347629556556754
$ws.Range("I8").Value = 0.3808887290954196
$ws.Range("J8").Value = 0.3808887290954196
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 476.6819177695423
$ws.Range("R8").Value = 4290.137259925881
$ws.Range("S8").Value = 0.03761585262981514
$ws.Range("T8").Value = 0.03761585262981515
$ws.Range("I9").Value = 0.3808887290954196
$ws.Range("J9").Value = 0.3808887290954196
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 1343.587732981265
$ws.Range("R9").Value = 12092.28959683139
$ws.Range("S9").Value = 0.1060249954425269
$ws.Range("T9").Value = 0.1060249954425269
$ws.Range("G10").Value = 7.213061
$ws.Range("H10").Value = 21.639183
$ws.Range("I10").Value = 0.1674173200514808
$ws.Range("J10").Value = 0.1674173200514808
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 13.84118843561967
$ws.Range("R10").Value = 124.570695920577
$ws.Range("S10").Value = 0.001092233804151733
$ws.Range("T10").Value = 0.001092233804151733
$ws.Range("G11").Value = 7.213061
$ws.Range("H11").Value = 21.639183
$ws.Range("I11").Value = 0.1674173200514808
$ws.Range("J11").Value = 0.1674173200514808
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 1307.644225279205
$ws.Range("R11").Value = 11768.79802751284
$ws.Range("S11").Value = 0.1031886267062305
$ws.Range("T11").Value = 0.1031886267062305
$ws.Range("G12").Value = 7.213061
$ws.Range("H12").Value = 21.639183
$ws.Range("I12").Value = 0.1674173200514808
$ws.Range("J12").Value = 0.1674173200514808
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 209.5226324483457
$ws.Range("R12").Value = 1885.703692035111
$ws.Range("S12").Value = 0.0165338188233904
$ws.Range("T12").Value = 0.0165338188233904
$ws.Range("G13").Value = 7.213061
$ws.Range("H13").Value = 21.639183
$ws.Range("I13").Value = 0.1674173200514808
$ws.Range("J13").Value = 0.1674173200514808
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 590.5658013141586
$ws.Range("R13").Value = 5315.092211827427
$ws.Range("S13").Value = 0.04660264071770815
$ws.Range("T13").Value = 0.04660264071770815
$ws.Range("G14").Value = 10.03371566666667
$ws.Range("H14").Value = 30.101147
$ws.Range("I14").Value = 0.2328855650980756
$ws.Range("J14").Value = 0.2328855650980756
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 19.25376053963256
$ws.Range("R14").Value = 173.283844856693
$ws.Range("S14").Value = 0.001519349889371541
$ws.Range("T14").Value = 0.001519349889371541
$ws.Range("G15").Value = 10.03371566666667
$ws.Range("H15").Value = 30.101147
$ws.Range("I15").Value = 0.2328855650980756
$ws.Range("J15").Value = 0.2328855650980756
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 1818.996172306065
$ws.Range("R15").Value = 16370.96555075458
$ws.Range("S15").Value = 0.1435403555306302
$ws.Range("T15").Value = 0.1435403555306302
$ws.Range("G16").Value = 10.03371566666667
$ws.Range("H16").Value = 30.101147
$ws.Range("I16").Value = 0.2328855650980756
$ws.Range("J16").Value = 0.2328855650980756
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 291.4560849711666
$ws.Range("R16").Value = 2623.104764740499
$ws.Range("S16").Value = 0.02299933924835524
$ws.Range("T16").Value = 0.02299933924835525
$ws.Range("G17").Value = 10.03371566666667
$ws.Range("H17").Value = 30.101147
$ws.Range("I17").Value = 0.2328855650980756
$ws.Range("J17").Value = 0.2328855650980756
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 821.5055068636502
$ws.Range("R17").Value = 7393.549561772853
$ws.Range("S17").Value = 0.06482652042971856
$ws.Range("T17").Value = 0.06482652042971856
